$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 64
$ws.Range("H64").Value = 5455.1665
$ws.Range("I64").Value = 3887.4285
$ws.Range("J64").Value = 7650
$ws.Range("K64").Value = 3887.4285
$ws.Range("L64").Value = 7650
$ws.Range("M64").Value = -3639.4285
$ws.Range("N64").Value = -8146
# row 67
$ws.Range("H67").Value = 5455.1665
$ws.Range("I67").Value = 3887.4285
$ws.Range("J67").Value = 7650
$ws.Range("K67").Value = 3887.4285
$ws.Range("L67").Value = 7650
$ws.Range("M67").Value = -3029.4285
$ws.Range("N67").Value = -9366
# row 111
$ws.Range("H111").Value = 3833
$ws.Range("I111").Value = 7500
$ws.Range("J111").Value = 1999.5
$ws.Range("K111").Value = 22500
$ws.Range("L111").Value = 5998.5
$ws.Range("M111").Value = -19433
$ws.Range("N111").Value = -12132.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1318.4482
$ws.Range("I2").Value = 903.95654
$ws.Range("K2").Value = 903.95654
$ws.Range("M2").Value = -790.95654
# row 32
$ws.Range("H32").Value = 224608.28
$ws.Range("I32").Value = 286811.03
$ws.Range("K32").Value = 286811.03
$ws.Range("M32").Value = -286524.03
# row 45
$ws.Range("H45").Value = 88013.25
$ws.Range("I45").Value = 146269.42
$ws.Range("J45").Value = 6454.6
$ws.Range("K45").Value = 146269.42
$ws.Range("L45").Value = 6454.6
$ws.Range("M45").Value = -145892.42
$ws.Range("N45").Value = -7208.6
# row 88
$ws.Range("H88").Value = 2264.5833
$ws.Range("J88").Value = 2176
$ws.Range("L88").Value = 2176
$ws.Range("N88").Value = -2988
# row 91
$ws.Range("H91").Value = 2264.5833
$ws.Range("J91").Value = 2176
$ws.Range("L91").Value = 2176
$ws.Range("N91").Value = -4984
# row 102
$ws.Range("H102").Value = 2064.7058
$ws.Range("I102").Value = 1952
$ws.Range("J102").Value = 2590.6667
$ws.Range("K102").Value = 1952
$ws.Range("L102").Value = 2590.6667
$ws.Range("M102").Value = -330
$ws.Range("N102").Value = -5834.6667
# row 116
$ws.Range("H116").Value = 1318.4482
$ws.Range("I116").Value = 903.95654
$ws.Range("K116").Value = 903.95654
$ws.Range("M116").Value = 1390.04346
# row 122
$ws.Range("H122").Value = 1496.3334
$ws.Range("I122").Value = 1496.2
$ws.Range("K122").Value = 4488.6
$ws.Range("M122").Value = -2038.6
# row 124
$ws.Range("H124").Value = 12500
$ws.Range("J124").Value = 12500
$ws.Range("L124").Value = 12500
$ws.Range("N124").Value = -22320
# row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1318.4482
$ws.Range("I3").Value = 903.95654
$ws.Range("K3").Value = 903.95654
$ws.Range("M3").Value = -789.95654
# row 107
$ws.Range("H107").Value = 6503.8477
$ws.Range("I107").Value = 7467.081
$ws.Range("J107").Value = 2543.889
$ws.Range("K107").Value = 7467.081
$ws.Range("L107").Value = 2543.889
$ws.Range("M107").Value = -5547.081
$ws.Range("N107").Value = -6383.889

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2883.75
$ws.Range("I31").Value = 2445.9355
$ws.Range("K31").Value = 2445.9355
$ws.Range("M31").Value = -2150.9355
# row 34
$ws.Range("H34").Value = 2883.75
$ws.Range("I34").Value = 2445.9355
$ws.Range("K34").Value = 2445.9355
$ws.Range("M34").Value = -2243.9355
# row 107
$ws.Range("H107").Value = 2846.2
$ws.Range("I107").Value = 2798.8333
$ws.Range("K107").Value = 2798.8333
$ws.Range("M107").Value = -878.8332999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 2039682.4
$ws.Range("I4").Value = 2396009.2
$ws.Range("J4").Value = 400579
$ws.Range("K4").Value = 7188027.600000001
$ws.Range("L4").Value = 1201737
$ws.Range("M4").Value = -7187915.600000001
$ws.Range("N4").Value = -1201961
# row 9
$ws.Range("H9").Value = 25112762
$ws.Range("I9").Value = 733.3333
$ws.Range("J9").Value = 37668776
$ws.Range("K9").Value = 2199.9999
$ws.Range("L9").Value = 113006328
$ws.Range("M9").Value = -1975.9999
$ws.Range("N9").Value = -113006776
# row 29
$ws.Range("H29").Value = 583.3333
$ws.Range("I29").Value = 375
$ws.Range("K29").Value = 1125
$ws.Range("M29").Value = -848
# row 52
$ws.Range("H52").Value = 1994.75
$ws.Range("J52").Value = 1994.75
$ws.Range("L52").Value = 5984.25
$ws.Range("N52").Value = -6516.25
# row 61
$ws.Range("H61").Value = 30068.75
$ws.Range("I61").Value = 99990
$ws.Range("J61").Value = 6761.6665
$ws.Range("K61").Value = 299970
$ws.Range("L61").Value = 20284.9995
$ws.Range("M61").Value = -299755
$ws.Range("N61").Value = -20714.9995
# row 140
$ws.Range("H140").Value = 11496125
$ws.Range("I140").Value = 13890318
$ws.Range("K140").Value = 41670954
$ws.Range("M140").Value = -41665774

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 62640376
$ws.Range("I80").Value = 222999.4
$ws.Range("J80").Value = 166669330
$ws.Range("K80").Value = 222999.4
$ws.Range("L80").Value = 166669330
$ws.Range("M80").Value = -222001.4
$ws.Range("N80").Value = -166671326
# row 83
$ws.Range("H83").Value = 62640376
$ws.Range("I83").Value = 222999.4
$ws.Range("J83").Value = 166669330
$ws.Range("K83").Value = 1114997
$ws.Range("L83").Value = 833346650
$ws.Range("M83").Value = -1110005
$ws.Range("N83").Value = -833356634
# row 101
$ws.Range("H101").Value = 108513.336
$ws.Range("J101").Value = 108513.336
$ws.Range("L101").Value = 108513.336
$ws.Range("N101").Value = -115003.336
# row 123
$ws.Range("H123").Value = 62500
$ws.Range("J123").Value = 62500
$ws.Range("L123").Value = 62500
$ws.Range("N123").Value = -67400
# row 126
$ws.Range("H126").Value = 2493.8
$ws.Range("I126").Value = 2493.8
$ws.Range("K126").Value = 7481.400000000001
$ws.Range("M126").Value = -5011.400000000001
# row 132
$ws.Range("H132").Value = 833533.9
$ws.Range("I132").Value = 8179.4736
$ws.Range("J132").Value = 2140345
$ws.Range("K132").Value = 24538.4208
$ws.Range("L132").Value = 6421035
$ws.Range("M132").Value = -22008.4208
$ws.Range("N132").Value = -6426095

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 2
$ws.Range("H2").Value = 13111
$ws.Range("J2").Value = 13111
$ws.Range("L2").Value = 13111
$ws.Range("N2").Value = -13335
# row 16
$ws.Range("H16").Value = 1145.7059
$ws.Range("I16").Value = 978.4666999999999
$ws.Range("K16").Value = 978.4666999999999
$ws.Range("M16").Value = -808.4666999999999
# row 101
$ws.Range("H101").Value = 23498
$ws.Range("J101").Value = 23498
$ws.Range("L101").Value = 23498
$ws.Range("N101").Value = -29988

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 51
$ws.Range("H51").Value = 2713.1428
$ws.Range("I51").Value = 2713.1428
$ws.Range("K51").Value = 2713.1428
$ws.Range("M51").Value = -2203.1428
# row 80
$ws.Range("H80").Value = 31150.5
$ws.Range("J80").Value = 31150.5
$ws.Range("L80").Value = 31150.5
$ws.Range("N80").Value = -33146.5
# row 83
$ws.Range("H83").Value = 31150.5
$ws.Range("J83").Value = 31150.5
$ws.Range("L83").Value = 93451.5
$ws.Range("N83").Value = -103435.5
# row 107
$ws.Range("H107").Value = 1299297
$ws.Range("I107").Value = 610.8421
$ws.Range("K107").Value = 1832.5263
$ws.Range("M107").Value = 87.47370000000001
# row 119
$ws.Range("H119").Value = 38665.668
$ws.Range("J119").Value = 38665.668
$ws.Range("L119").Value = 38665.668
$ws.Range("N119").Value = -48341.668
# row 132
$ws.Range("H132").Value = 2242.353
$ws.Range("I132").Value = 1851.7838
$ws.Range("K132").Value = 5555.3514
$ws.Range("M132").Value = -3025.3514
